# ck_retailers.xlsx - "added fields to retailers, map and print work done."
#
# Business-data change: the last retailer row (row 5, "Kesav Electronics")
# gets a corrected retailer_code and a new dse_code assignment.
#   A5 (retailer_code): "KG10"     -> "KG11"
#   C5 (dse_code)      : "sunil006" -> "sunil009"
# Everything else in the sheet (headers, rows 2-4, route_no) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "KG11"
$ws.Range("C5").Value = "sunil009"

# Leave the cursor/selection on the row that was edited, as recorded in the
# saved view state (sheetView/selection moved from A6 to A5).
$ws.Range("A5").Select()
